# Run job commit message
#
# The workbook records JOB_START / JOB_END timestamps for each task of a
# "DAMAGE_RETURN" job run. This commit represents a fresh run of the job:
# every task that previously completed (had both a start and end time) gets
# new start/end timestamps from the new run, and the task that was in
# progress when the sheet was captured (TASK_NO 24, row 25) gets only a new
# JOB_START with JOB_END left blank (it had not finished yet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# TASK_NO 3 (row 4) - use warehouse demo_WH;
$ws.Range("H4").Value = "2021-09-02 09:42:06.992"
$ws.Range("I4").Value = "2021-09-02 09:42:08.805"

# TASK_NO 4 (row 5) - use database Demo1_DB;
$ws.Range("H5").Value = "2021-09-02 09:42:09.228"
$ws.Range("I5").Value = "2021-09-02 09:42:10.947"

# TASK_NO 6 (row 7) - create or replace table ETL_SALES_TABLE_1 ...
$ws.Range("H7").Value = "2021-09-02 09:42:11.367"
$ws.Range("I7").Value = "2021-09-02 09:42:13.491"

# TASK_NO 7 (row 8) - Insert Into ETL_SALES_TABLE_1 ...
$ws.Range("H8").Value = "2021-09-02 09:42:13.872"
$ws.Range("I8").Value = "2021-09-02 09:42:14.847"

# TASK_NO 10 (row 11) - copy into ETL_SALES_TABLE_1 ...
$ws.Range("H11").Value = "2021-09-02 09:42:15.132"
$ws.Range("I11").Value = "2021-09-02 09:42:18.079"

# TASK_NO 13 (row 14) - CREATE OR REPLACE TABLE ETL_SALES_TABLE_2 ...
$ws.Range("H14").Value = "2021-09-02 09:42:18.443"
$ws.Range("I14").Value = "2021-09-02 09:42:20.197"

# TASK_NO 14 (row 15) - Insert Into ETL_SALES_TABLE_2 ...
$ws.Range("H15").Value = "2021-09-02 09:42:20.473"
$ws.Range("I15").Value = "2021-09-02 09:42:21.365"

# TASK_NO 15 (row 16) - insert into ETL_SALES_TABLE_2 select * ...
$ws.Range("H16").Value = "2021-09-02 09:42:21.635"
$ws.Range("I16").Value = "2021-09-02 09:42:23.068"

# TASK_NO 17 (row 18) - CREATE OR REPLACE TABLE ETL_SALES_TABLE_3 CLONE ...
$ws.Range("H18").Value = "2021-09-02 09:42:23.374"
$ws.Range("I18").Value = "2021-09-02 09:42:25.889"

# TASK_NO 18 (row 19) - Update ETL_SALES_TABLE_3 set Region= ...
$ws.Range("H19").Value = "2021-09-02 09:42:26.343"
$ws.Range("I19").Value = "2021-09-02 09:42:27.42"

# TASK_NO 19 (row 20) - Update ETL_SALES_TABLE_3 set SALE_AMOUNT= ...
$ws.Range("H20").Value = "2021-09-02 09:42:27.831"
$ws.Range("I20").Value = "2021-09-02 09:42:28.882"

# TASK_NO 20 (row 21) - CREATE OR REPLACE TABLE ETL_SALES_TABLE_4 CLONE ...
$ws.Range("H21").Value = "2021-09-02 09:42:29.265"
$ws.Range("I21").Value = "2021-09-02 09:42:31.785"

# TASK_NO 21 (row 22) - Update ETL_SALES_TABLE_4 set Region= ...
$ws.Range("H22").Value = "2021-09-02 09:42:32.152"
$ws.Range("I22").Value = "2021-09-02 09:42:33.169"

# TASK_NO 22 (row 23) - Update ETL_SALES_TABLE_4 set RETURN_AMOUNT= ...
$ws.Range("H23").Value = "2021-09-02 09:42:33.531"
$ws.Range("I23").Value = "2021-09-02 09:42:34.466"

# TASK_NO 24 (row 25) - CREATE OR REPLACE TABLE ETL_SALES_TABLE_4 CLONE ...
# This task only started in the new run; it has not completed, so JOB_END
# is left blank.
$ws.Range("H25").Value = "2021-09-02 09:42:34.819"
$ws.Range("I25").Value = ""
